$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.014.21"
$ws.Range("E2").Value = "  -1.13%  "

$ws.Range("D3").Value = "'1.885.73"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "'330.79"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "'0.4600"
$ws.Range("E7").Value = "  -3.48%  "

$ws.Range("D8").Value = "'0.4069"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").Value = "'47.76"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").Value = "'0.07980"
$ws.Range("E10").Value = "  -2.71%  "

$ws.Range("D11").Value = "'0.9896"
$ws.Range("E11").Value = "  -4.73%  "

$ws.Range("E12").Value = "  -4.35%  "

$ws.Range("D13").Value = "'1.896.22"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "'5.915"
$ws.Range("E14").Value = "  -3.19%  "

$ws.Range("D15").Value = "'7.067"
$ws.Range("E15").Value = "  -4.60%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "'88.40"
$ws.Range("E17").Value = "  -4.02%  "

$ws.Range("D18").Value = "'0.00001031"
$ws.Range("E18").Value = "  -2.70%  "

$ws.Range("D19").Value = "'0.06574"
$ws.Range("E19").Value = "  -1.11%  "

$ws.Range("D20").Value = "'17.45"
$ws.Range("E20").Value = "  -2.90%  "

$ws.Range("D21").Value = "'0.9989"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").Value = "'29.040.38"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").Value = "'5.417"
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("D24").Value = "'11.49"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("D26").Value = "'2.094.28"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").Value = "'156.71"
$ws.Range("E27").Value = "  -2.64%  "

$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("D29").Value = "'2.084"
$ws.Range("E29").Value = "  -5.42%  "

$ws.Range("D30").Value = "'5.481"
$ws.Range("E30").Value = "  -1.82%  "

$ws.Range("D31").Value = "'117.56"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("D32").Value = "'1.009"
$ws.Range("E32").Value = "  -1.57%  "

$ws.Range("D33").Value = "'0.09325"
$ws.Range("E33").Value = "  -2.75%  "

$ws.Range("D34").Value = "'3.603"
$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").Value = "'1.404"
$ws.Range("E35").Value = "  -3.14%  "

$ws.Range("D36").Value = "'5.270"
$ws.Range("E36").Value = "  -3.50%  "

$ws.Range("D37").Value = "'0.06061"
$ws.Range("E37").Value = "  -2.65%  "

$ws.Range("D38").Value = "'0.02223"
$ws.Range("E38").Value = "  -3.34%  "

$ws.Range("D39").Value = "'8.295"
$ws.Range("E39").Value = "  -4.64%  "

$ws.Range("D40").Value = "'1.172"
$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").Value = "'0.5772"
$ws.Range("E42").Value = "  -5.02%  "

$ws.Range("D43").Value = "'0.1826"
$ws.Range("E43").Value = "  -4.39%  "

$ws.Range("D44").Value = "'10.08"
$ws.Range("E44").Value = "  -5.86%  "

$ws.Range("D45").Value = "'1.249"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").Value = "'0.07500"
$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.261"
$ws.Range("E47").Value = "  -3.30%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'12.04"
$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("D49").Value = "'0.5450"
$ws.Range("E49").Value = "  -3.97%  "

$ws.Range("D50").Value = "'1.902"
$ws.Range("E50").Value = "  -4.53%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'111.04"
$ws.Range("E51").Value = "  -2.17%  "
